# Fill in the missing xG_home / xG_away / goals_home / goals_away values
# for the six most recent Fiorentina matches (rows 10-15), which previously
# only had home/away team names filled in.
#
# All of D/E/F/G in this sheet are stored as text (shared strings), even
# the goal counts that look numeric, so each value is written via a
# helper cell that is temporarily formatted as Text (to stop Excel from
# auto-converting the numeric-looking text back into a Number), then
# copied as values-only onto the real target cell so the target cell's
# style stays untouched (matching the existing un-styled D2:G9 cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("Z1")

function Set-TextValue($cellRef, $text) {
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

$rows = @(
    @{ Row = 10; D = "2.85768";  E = "0.602734"; F = "2"; G = "0" },
    @{ Row = 11; D = "2.79212";  E = "1.37691";  F = "1"; G = "1" },
    @{ Row = 12; D = "2.8038";   E = "0.461794"; F = "3"; G = "0" },
    @{ Row = 13; D = "1.71835";  E = "1.7103";   F = "1"; G = "1" },
    @{ Row = 14; D = "1.68253";  E = "1.24682";  F = "1"; G = "1" },
    @{ Row = 15; D = "0.388088"; E = "1.53117";  F = "0"; G = "3" }
)

foreach ($r in $rows) {
    Set-TextValue ("D" + $r.Row) $r.D
    Set-TextValue ("E" + $r.Row) $r.E
    Set-TextValue ("F" + $r.Row) $r.F
    Set-TextValue ("G" + $r.Row) $r.G
}

$helper.Clear() | Out-Null
$excel.CutCopyMode = 0
